# Update "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, as published at the new data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1496
$ws.Range("F5").Value  = 7527
$ws.Range("F6").Value  = 80
$ws.Range("F7").Value  = 4801
$ws.Range("F11").Value = 1498
$ws.Range("F22").Value = 1173
$ws.Range("F26").Value = 1232
$ws.Range("F31").Value = 188
$ws.Range("F45").Value = 21

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 1733
$ws.Range("F45").Value = 77

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 680
$ws.Range("F9").Value = 2526

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 1496
$ws.Range("F7").Value  = 680
$ws.Range("F8").Value  = 680
$ws.Range("F9").Value  = 7527
$ws.Range("F10").Value = 80
$ws.Range("F11").Value = 4801
$ws.Range("F15").Value = 1498
$ws.Range("F19").Value = 2526
$ws.Range("F26").Value = 1173
$ws.Range("F30").Value = 1232
$ws.Range("F32").Value = 188
$ws.Range("F49").Value = 77
